$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each touched column-D cell to Text format before writing its new
# value, so numeric-looking strings (e.g. "0.999", "1.00", "123.22") are
# preserved exactly as text instead of being parsed into floating point
# numbers by Excel (which would lose trailing zeros / exact formatting).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '50.000.49'
$ws.Range("E2").Value = '  +3.95%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.648.71'
$ws.Range("E3").Value = '  +6.30%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '114.19'
$ws.Range("E5").Value = '  +8.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '326.50'
$ws.Range("E6").Value = '  +3.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.529'
$ws.Range("E7").Value = '  +2.24%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  +4.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.13'
$ws.Range("E10").Value = '  +6.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.17'
$ws.Range("E11").Value = '  -0.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0825'
$ws.Range("E12").Value = '  +2.97%  '
$ws.Range("E13").Value = '  +0.94%  '
$ws.Range("E14").Value = '  +4.83%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.061.86'
$ws.Range("E15").Value = '  +6.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.671.46'
$ws.Range("E16").Value = '  +7.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.873'
$ws.Range("E17").Value = '  +5.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '49.896.02'
$ws.Range("E18").Value = '  +4.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.22'
$ws.Range("E19").Value = '  +3.29%  '
$ws.Range("E20").Value = '  +3.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.92'
$ws.Range("E21").Value = '  -1.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0958'
$ws.Range("E22").Value = '  +3.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.10'
$ws.Range("E23").Value = '  +1.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '276.42'
$ws.Range("E24").Value = '  +2.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.58'
$ws.Range("E25").Value = '  +3.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.84'
$ws.Range("E26").Value = '  +4.34%  '
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.04'
$ws.Range("E28").Value = '  +3.61%  '
$ws.Range("E29").Value = '  +1.50%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.23'
$ws.Range("E30").Value = '  +5.27%  '
$ws.Range("E31").Value = '  +2.86%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.35'
$ws.Range("E32").Value = '  +2.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.47'
$ws.Range("E33").Value = '  +3.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.60'
$ws.Range("E34").Value = '  +3.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0809'
$ws.Range("E35").Value = '  +5.12%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.08'
$ws.Range("E36").Value = '  +11.44%  '
$ws.Range("B37").Value = 'FirstDigitalUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("E38").Value = '  +7.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.11'
$ws.Range("E39").Value = '  +8.50%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.113'
$ws.Range("E40").Value = '  +2.19%  '
$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '123.22'
$ws.Range("E41").Value = '  +0.58%  '
$ws.Range("E42").Value = '  +0.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.02'
$ws.Range("E43").Value = '  -1.83%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0316'
$ws.Range("E44").Value = '  +5.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.086.29'
$ws.Range("E45").Value = '  +4.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.33'
$ws.Range("E46").Value = '  +6.78%  '
$ws.Range("E47").Value = '  +15.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.99'
$ws.Range("E48").Value = '  +5.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.14'
$ws.Range("E49").Value = '  +2.93%  '
$ws.Range("E50").Value = '  +5.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '59.61'
$ws.Range("E51").Value = '  +6.39%  '
